# Update the "テスト仕様書" (test specification) sheet: for test case rows
# 5 through 65, mark retest-needed ("再試要否") as "否" (No), stamp the
# implementation date ("実施日") and record the approver ("承認者").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Implementation date (serial 46045 == 2026-01-23), stored without a time
# component so the cell keeps a clean integer date serial.
$implementationDate = (Get-Date -Year 2026 -Month 1 -Day 23 -Hour 0 -Minute 0 -Second 0).Date

for ($row = 5; $row -le 65; $row++) {
    $ws.Range("I$row").Value = "否"
    $ws.Range("J$row").Value = $implementationDate
    $ws.Range("L$row").Value = "清山"
}
